# Auto-generated edit script: updates crypto price/volume figures and
# re-ranks several coin rows to match the refreshed "Updated symbol list"
# snapshot (commit: "Updated symbol list on Sun Jan  1 17:17:09 UTC 2023").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay plain Text (matching the source inline-string
    # cells) instead of letting Excel auto-coerce numeric- or percent-
    # looking strings ("244.91", "-0.81%", "3,762.43%", ...) into numbers.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '244.91'
Set-TextValue $ws.Range('E2') '-0.81%'

Set-TextValue $ws.Range('D3') '27.05'
Set-TextValue $ws.Range('E3') '1.87%'

Set-TextValue $ws.Range('D4') '5.055'
Set-TextValue $ws.Range('E4') '-0.56%'

Set-TextValue $ws.Range('D5') '0.05689'
Set-TextValue $ws.Range('E5') '1.34%'

Set-TextValue $ws.Range('D6') '6.471'
Set-TextValue $ws.Range('E6') '-0.19%'

Set-TextValue $ws.Range('D7') '0.8205'
Set-TextValue $ws.Range('E7') '0.91%'

Set-TextValue $ws.Range('D8') '0.8384'
Set-TextValue $ws.Range('E8') '-0.69%'

Set-TextValue $ws.Range('D9') '0.1327'
Set-TextValue $ws.Range('E9') '-1.14%'

Set-TextValue $ws.Range('D10') '0.06916'
Set-TextValue $ws.Range('E10') '-0.78%'

Set-TextValue $ws.Range('D11') '0.02857'
Set-TextValue $ws.Range('E11') '-0.43%'

Set-TextValue $ws.Range('D12') '0.09396'
Set-TextValue $ws.Range('E12') '-0.13%'

Set-TextValue $ws.Range('D13') '0.001525'
Set-TextValue $ws.Range('E13') '0.97%'

Set-TextValue $ws.Range('D14') '0.04078'
Set-TextValue $ws.Range('E14') '-12.65%'

Set-TextValue $ws.Range('D15') '0.0005987'
Set-TextValue $ws.Range('E15') '-0.23%'

Set-TextValue $ws.Range('D16') '0.006189'
Set-TextValue $ws.Range('E16') '0.86%'

$ws.Range('B17').Value = 'UpBots'
$ws.Range('C17').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue $ws.Range('D17') '0.007486'
Set-TextValue $ws.Range('E17') '3,762.43%'

$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D18') '3.509'
Set-TextValue $ws.Range('E18') '-2.17%'

$ws.Range('B19').Value = 'GateToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D19') '3.005'
Set-TextValue $ws.Range('E19') '-0.20%'

$ws.Range('B20').Value = 'BTSEToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D20') '2.228'
Set-TextValue $ws.Range('E20') '5.17%'

$ws.Range('B21').Value = 'BitpandaEcosystemToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D21') '0.3177'
Set-TextValue $ws.Range('E21') '0.65%'

$ws.Range('B22').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C22').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D22') '0.03186'
Set-TextValue $ws.Range('E22') '-0.45%'

$ws.Range('B23').Value = 'ProBitToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D23') '0.1296'
Set-TextValue $ws.Range('E23') '-1.83%'

$ws.Range('B24').Value = 'MCDex'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D24') '3.582'
Set-TextValue $ws.Range('E24') '-4.30%'

$ws.Range('B25').Value = 'ZBToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range('D25') '0.1373'
Set-TextValue $ws.Range('E25') '1.72%'

$ws.Range('B26').Value = 'BitKan'
$ws.Range('C26').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D26') '0.001217'
Set-TextValue $ws.Range('E26') '-2.83%'

$ws.Range('B27').Value = 'HotbitToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D27') '0.003957'
Set-TextValue $ws.Range('E27') '-13.95%'

$ws.Range('B28').Value = 'NitroEx'
$ws.Range('C28').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws.Range('D28') '0.00009795'
Set-TextValue $ws.Range('E28') '1.97%'

Set-TextValue $ws.Range('D40') '0.03687'
Set-TextValue $ws.Range('E40') '0.27%'

$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range('D41') '0.1055'
Set-TextValue $ws.Range('E41') '-0.32%'

$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range('D42') '0.005987'
Set-TextValue $ws.Range('E42') '-3.34%'

Set-TextValue $ws.Range('D43') '0.002334'
Set-TextValue $ws.Range('E43') '-10.35%'

Set-TextValue $ws.Range('D44') '0.009372'
Set-TextValue $ws.Range('E44') '5.19%'

Set-TextValue $ws.Range('D45') '0.00005211'
Set-TextValue $ws.Range('E45') '-1.56%'

Set-TextValue $ws.Range('E46') '-0.11%'

Set-TextValue $ws.Range('E47') '-15.46%'

Set-TextValue $ws.Range('D48') '0.002596'
Set-TextValue $ws.Range('E48') '2.99%'

Set-TextValue $ws.Range('E49') '-0.11%'

Set-TextValue $ws.Range('E50') '-0.11%'
